$d = $word.ActiveDocument

# The report-submission date line reads "« 11 »   декабря  2021 г." — the
# day-of-month run ("11") is updated to "18" (the report was actually
# submitted on the 18th, not the 11th). It is the sole run in the whole
# document whose text is exactly "11", so a whole-word Find/Replace over
# the document body targets it uniquely and keeps the surrounding
# characters/formatting (the « », underline, "декабря  2021 г.", etc.)
# untouched.
$d.Content.Find.Execute("11", $true, $false, $false, $false, $false, $true, 1, $false, "18", 2) | Out-Null
